$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new row at position 30 (shifts the old row 30+ content, incl. the
#    signature block at rows 34-35, down by one row to 35-36, and extends the
#    merged cell ranges/table automatically).
$ws.Rows("30").Insert()

# 2) Fix up the borders/style of the (now 15-row) data table: row 30 should
#    carry the "closing" bottom-border style that used to belong to row 29,
#    and row 29 should fall back to the regular "middle" row style (like row 28).
$ws.Range("B29:J29").Copy($ws.Range("B30:J30"))
$ws.Range("B28:J28").Copy($ws.Range("B29:J29"))

# 3) Rewrite the worker/period rows 16-30 with the updated, re-sorted dataset
#    (Nancy's periods 2407-2505 in order, then Ronal 2506, Nancy 2506, Nancy
#    2507, and the newly added Ronal 2508).
$tableRows = @(
    @("CC", "1094245611", "NANCY ASTRID TIBAMOZA PEÑA", "2407", 177709, 4442734),
    @("CC", "1094245611", "NANCY ASTRID TIBAMOZA PEÑA", "2408", 177709, 4442734),
    @("CC", "1094245611", "NANCY ASTRID TIBAMOZA PEÑA", "2409", 177709, 4442734),
    @("CC", "1094245611", "NANCY ASTRID TIBAMOZA PEÑA", "2410", 177709, 4442734),
    @("CC", "1094245611", "NANCY ASTRID TIBAMOZA PEÑA", "2411", 177709, 4442734),
    @("CC", "1094245611", "NANCY ASTRID TIBAMOZA PEÑA", "2412", 177709, 4442734),
    @("CC", "1094245611", "NANCY ASTRID TIBAMOZA PEÑA", "2501", 177709, 4442734),
    @("CC", "1094245611", "NANCY ASTRID TIBAMOZA PEÑA", "2502", 177709, 4442734),
    @("CC", "1094245611", "NANCY ASTRID TIBAMOZA PEÑA", "2503", 177709, 4442734),
    @("CC", "1094245611", "NANCY ASTRID TIBAMOZA PEÑA", "2504", 177709, 4442734),
    @("CC", "1094245611", "NANCY ASTRID TIBAMOZA PEÑA", "2505", 177709, 4442734),
    @("CC", "73199528",   "RONAL OROZCO CONTRERAS",     "2506", 68328,  1708200),
    @("CC", "1094245611", "NANCY ASTRID TIBAMOZA PEÑA", "2506", 177709, 4442734),
    @("CC", "1094245611", "NANCY ASTRID TIBAMOZA PEÑA", "2507", 177709, 4442734),
    @("CC", "73199528",   "RONAL OROZCO CONTRERAS",     "2508", 177709, 4442734)
)

$r = 16
foreach ($row in $tableRows) {
    $ws.Cells.Item($r, 2).Value = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
    $ws.Cells.Item($r, 6).Value = $row[4]
    $ws.Cells.Item($r, 7).Value = $row[5]
    $r = $r + 1
}

# 4) Update the summary header figures: total "Valor Mora" and period count.
$ws.Range("E11").Value = 2556254
$ws.Range("F13").Value = 14
